$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "17:15:07", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:16:17", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:16:27", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:16:38", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:16:48", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:16:59", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        # Column A holds a literal date-like string (e.g. "2026-02-01").
        # Force text storage so Excel doesn't reinterpret it as a date
        # serial number, then strip the resulting number-format style so
        # the cell keeps behaving like its plain-text neighbours.
        if ($c -eq 0) {
            $cell.NumberFormat = "@"
            $cell.Value = $data[$c]
            $cell.ClearFormats()
        } else {
            $cell.Value = $data[$c]
        }
    }
}
